$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.265887379646301
$ws.Range("B1").Value = 3.475752115249634
$ws.Range("C1").Value = 4.56719970703125
$ws.Range("D1").Value = 2.623500823974609
$ws.Range("E1").Value = 2.203383922576904
